$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Medico Clinico" (row 6) and "Intensivista" (row 4) rows,
# shifting the remaining rows up.
$ws.Range("A6:B6").EntireRow.Delete()
$ws.Range("A4:B4").EntireRow.Delete()

# Reorder the remaining specialties (Cirujano before Enfermera) and add the
# "Cantidad" counts in column B for each.
$ws.Range("A2").Value = "Cirujano"
$ws.Range("B2").Value = 0
$ws.Range("A3").Value = "Enfermera"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "Odontologa"
$ws.Range("B4").Value = 0
